$wb = $excel.ActiveWorkbook

# --- Sheet "Schedule": update computed Cost ($) and Unit Cost ($/ML) for rows 3-4 ---
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("E3").Value = -32.464926
$wsSchedule.Range("F3").Value = -1.073575595238095
$wsSchedule.Range("E4").Value = 466.16691225
$wsSchedule.Range("F4").Value = 30.83114499007937

# --- Sheet "Detailed": refresh Price column (and Type for rows that flipped from forecast to historical) ---
$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B37").Value = 24.43616
$wsDetailed.Range("B38").Value = 21.27683
$wsDetailed.Range("B39").Value = 67.18411999999999
$wsDetailed.Range("C39").Value = "historical"
$wsDetailed.Range("B40").Value = 105.56944
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 109.41
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 117.58107
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 108.89
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 84.79000000000001
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 84.79000000000001
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("B46").Value = 84.79000000000001
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 123.97173
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 138.42
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 105.79
$wsDetailed.Range("B50").Value = 97.93371
$wsDetailed.Range("B51").Value = 104.37914
$wsDetailed.Range("B52").Value = 96.02821
$wsDetailed.Range("B53").Value = 103.99643
$wsDetailed.Range("B54").Value = 85.07574
$wsDetailed.Range("B55").Value = 84.79000000000001
$wsDetailed.Range("B56").Value = 84.79000000000001
$wsDetailed.Range("B57").Value = 84.79000000000001
$wsDetailed.Range("B58").Value = 84.79000000000001
$wsDetailed.Range("B59").Value = 84.69816
$wsDetailed.Range("B60").Value = 84.79000000000001
$wsDetailed.Range("B61").Value = 105.79
$wsDetailed.Range("B62").Value = 108.89
$wsDetailed.Range("B63").Value = 101.43757
$wsDetailed.Range("B64").Value = 57.06015
$wsDetailed.Range("B65").Value = 36.05919
$wsDetailed.Range("B66").Value = 8.62425
$wsDetailed.Range("B67").Value = 8.43248
$wsDetailed.Range("B68").Value = 0.7
$wsDetailed.Range("B70").Value = 0.01116
$wsDetailed.Range("B71").Value = -5.50985
$wsDetailed.Range("B72").Value = -5.50985
$wsDetailed.Range("B73").Value = -5.47284
$wsDetailed.Range("B74").Value = -5.50985
$wsDetailed.Range("B75").Value = -5.33373
$wsDetailed.Range("B76").Value = -1.15634
$wsDetailed.Range("B77").Value = -3.6481
$wsDetailed.Range("B79").Value = -2.83895
$wsDetailed.Range("B80").Value = -5.08589
$wsDetailed.Range("B81").Value = -6
$wsDetailed.Range("B82").Value = -5.45139
$wsDetailed.Range("B83").Value = -4.80816
$wsDetailed.Range("B84").Value = 13.81347
$wsDetailed.Range("B85").Value = 32.72192
$wsDetailed.Range("B86").Value = 59.45881
$wsDetailed.Range("B87").Value = 62.71208
$wsDetailed.Range("B88").Value = 69.06494000000001
$wsDetailed.Range("B89").Value = 73.19
$wsDetailed.Range("B90").Value = 73.2
$wsDetailed.Range("B91").Value = 64.89
$wsDetailed.Range("B92").Value = 62.92943
$wsDetailed.Range("B93").Value = 62.60589
$wsDetailed.Range("B94").Value = 57.06
$wsDetailed.Range("B95").Value = 49.60785
$wsDetailed.Range("B96").Value = 50.76674
